$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6625
$ws.Range("F3").Value = 780
$ws.Range("F5").Value = 124
$ws.Range("F6").Value = 656
$ws.Range("F8").Value = 50
$ws.Range("F9").Value = 824
$ws.Range("F10").Value = 1284
$ws.Range("F14").Value = 519
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 359
$ws.Range("F17").Value = 1047
$ws.Range("F19").Value = 707
$ws.Range("F21").Value = 431
$ws.Range("F23").Value = 1108
$ws.Range("F24").Value = 218
$ws.Range("F25").Value = 2314
$ws.Range("F27").Value = 170
$ws.Range("F28").Value = 423
$ws.Range("F30").Value = 3749
$ws.Range("F32").Value = 688

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 35
$ws.Range("F6").Value = 736
$ws.Range("F11").Value = 137
$ws.Range("G11").Value = 180
$ws.Range("F16").Value = 6
$ws.Range("F24").Value = 220

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1230
$ws.Range("F5").Value = 1612
$ws.Range("F8").Value = 924

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1230
$ws.Range("F4").Value = 1612
$ws.Range("F7").Value = 924
$ws.Range("F9").Value = 6625
$ws.Range("F10").Value = 35
$ws.Range("F11").Value = 780
$ws.Range("F12").Value = 736
$ws.Range("F13").Value = 124
$ws.Range("F14").Value = 656
$ws.Range("F15").Value = 656
$ws.Range("F17").Value = 50
$ws.Range("F18").Value = 824
$ws.Range("F21").Value = 137
$ws.Range("G21").Value = 180
$ws.Range("F22").Value = 137
$ws.Range("G22").Value = 180
$ws.Range("F25").Value = 1284
$ws.Range("F29").Value = 519
$ws.Range("F30").Value = 6
$ws.Range("F32").Value = 359
$ws.Range("F33").Value = 1047
$ws.Range("F36").Value = 707
$ws.Range("F38").Value = 431
$ws.Range("F40").Value = 220
$ws.Range("F41").Value = 1108
$ws.Range("F42").Value = 218
$ws.Range("F43").Value = 2314
$ws.Range("F47").Value = 170
$ws.Range("F48").Value = 423
$ws.Range("F49").Value = 3749
$ws.Range("F52").Value = 688
